$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): switch from "duplicate of first data row" layout to a
#     proper field-name header, and extend with the new deposit/legislator columns ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Match the bold / centered / bordered header look used for B1:G1 on the new H1:M1 cells
$ws.Range("H1:M1").Font.Bold = $true
$ws.Range("H1:M1").HorizontalAlignment = -4108
$ws.Range("H1:M1").VerticalAlignment = -4160
$ws.Range("H1:M1").Borders.LineStyle = 1

# --- Data rows 2-8: bank / deposit_type / currency / owner / total stay (shifted),
#     and new property_category / category / date / legislator_name / legislator_id /
#     source_file / index columns are appended (same pattern as the other sheets) ---

# Row 2 -- 安泰商業銀行 / 活期儲蓄存款
$ws.Range("B2").Value = "安泰商業銀行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "羅淑蕾"
$ws.Range("F2").Value = 6098718
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2011-11-21"
$ws.Range("J2").Value = "羅淑蕾"
$ws.Range("K2").Value = 1638
$ws.Range("L2").Value = "tmpa0031"
$ws.Range("M2").Value = 90

# Row 3 -- 安泰商業銀行 / 支票存款
$ws.Range("B3").Value = "安泰商業銀行"
$ws.Range("C3").Value = "支票存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "羅淑蕾"
$ws.Range("F3").Value = 161
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2011-11-21"
$ws.Range("J3").Value = "羅淑蕾"
$ws.Range("K3").Value = 1638
$ws.Range("L3").Value = "tmpa0031"
$ws.Range("M3").Value = 91

# Row 4 -- 華泰商業銀行 / 活期儲蓄存款
$ws.Range("B4").Value = "華泰商業銀行"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "羅淑蕾"
$ws.Range("F4").Value = 29425
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2011-11-21"
$ws.Range("J4").Value = "羅淑蕾"
$ws.Range("K4").Value = 1638
$ws.Range("L4").Value = "tmpa0031"
$ws.Range("M4").Value = 92

# Row 5 -- 華泰商業銀行 / 支票存款
$ws.Range("B5").Value = "華泰商業銀行"
$ws.Range("C5").Value = "支票存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "羅淑蕾"
$ws.Range("F5").Value = 59831
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2011-11-21"
$ws.Range("J5").Value = "羅淑蕾"
$ws.Range("K5").Value = 1638
$ws.Range("L5").Value = "tmpa0031"
$ws.Range("M5").Value = 93

# Row 6 -- 中國銀行 / 綜合存款 (人民幣)
$ws.Range("B6").Value = "中國銀行"
$ws.Range("C6").Value = "綜合存款"
$ws.Range("D6").Value = "人民幣"
$ws.Range("E6").Value = "羅淑蕾"
$ws.Range("F6").Value = 7879500
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2011-11-21"
$ws.Range("J6").Value = "羅淑蕾"
$ws.Range("K6").Value = 1638
$ws.Range("L6").Value = "tmpa0031"
$ws.Range("M6").Value = 94

# Row 7 -- 美商美國銀行 / 綜合存款 (美金)
$ws.Range("B7").Value = "美商美國銀行"
$ws.Range("C7").Value = "綜合存款"
$ws.Range("D7").Value = "美金"
$ws.Range("E7").Value = "羅淑蕾"
$ws.Range("F7").Value = 30533705.8
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2011-11-21"
$ws.Range("J7").Value = "羅淑蕾"
$ws.Range("K7").Value = 1638
$ws.Range("L7").Value = "tmpa0031"
$ws.Range("M7").Value = 95

# Row 8 -- 永豐商業銀行敦南分行 / 活期儲蓄存款
$ws.Range("B8").Value = "永豐商業銀行敦南分行"
$ws.Range("C8").Value = "活期儲蓄存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "羅淑蕾"
$ws.Range("F8").Value = 1593775
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2011-11-21"
$ws.Range("J8").Value = "羅淑蕾"
$ws.Range("K8").Value = 1638
$ws.Range("L8").Value = "tmpa0031"
$ws.Range("M8").Value = 96
